$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109, shifting existing rows 109-119 down to 110-120
$ws.Rows("109:109").Insert()

# Populate the newly inserted row 109 with the new weekly data point
$ws.Range("A109").Value = 1
$ws.Range("B109").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C109").Value = "Arica y Parinacota"
$ws.Range("D109").Value = 45106
$ws.Range("E109").Value = 15
$ws.Range("F109").Value = 100112040
$ws.Range("G109").Value = "Cilantro"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 225
$ws.Range("K109").Value = 2000
$ws.Range("L109").Value = 2500
$ws.Range("M109").Value = 2278
$ws.Range("N109").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O109").Value = "Región de Arica y Parinacota"
$ws.Range("P109").Value = 1139
$ws.Range("Q109").Value = 2
$ws.Range("R109").Value = "Hortaliza"
